# "Back UP 14th May" - append duplicated/extra contact rows (Pawan/Awasthi
# series) to the NewContact sheet, rows 2-20.
#
# Columns are filled one at a time (all of B, then all of C, then A, then D)
# to match how the shared-string table was actually built (all "PawanN"
# strings first, then all "AwasthiN" strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewContact")

$titles = @("Dr.", "Miss", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.", "Mr.")
$firstNames = @("Pawan1", "Pawan2", "Pawan3", "Pawan4", "Pawan5", "Pawan6", "Pawan7", "Pawan8", "Pawan9", "Pawan10", "Pawan11", "Pawan12", "Pawan13", "Pawan14", "Pawan15", "Pawan16", "Pawan17", "Pawan18", "Pawan19")
$lastNames = @("Awasthi1", "Awasthi2", "Awasthi3", "Awasthi4", "Awasthi5", "Awasthi6", "Awasthi7", "Awasthi8", "Awasthi9", "Awasthi10", "Awasthi11", "Awasthi12", "Awasthi13", "Awasthi14", "Awasthi15", "Awasthi16", "Awasthi17", "Awasthi18", "Awasthi19")
$companies = @("HCL", "Google", "Ebay", "HCL", "Google", "Ebay", "HCL", "Google", "Ebay", "HCL", "Google", "Ebay", "HCL", "Google", "Ebay", "HCL", "Google", "Ebay", "Ebay")

for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $firstNames[$i]
}
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $lastNames[$i]
}
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $titles[$i]
}
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $companies[$i]
}

# Scroll the view so row 4 is at the top, then leave the active selection on G15
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G15").Select()
